# Commit extend report group
# Extends the "CheckTitleOfHomePageLink" report table on Sheet1 with a new
# "Medical, Health & Beauty" row, adds a second report group
# "CheckTitleOfHomePageLink2" below it, and adds a brand-new "Sheet2"
# worksheet containing a "CheckStaticLinks" report group.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet1: update existing "Health & Beauty" block --------------------
# Row 5 used to be "Metallurgy & Chemicals" with no Runmode flag; give it one.
$ws1.Range("A5").Value = "N"

# Row 4 used to be "Health & Beauty" / "healthcare_cosmetics_personal_care";
# turn it into the new "Medical, Health & Beauty" row with a Runmode flag.
$ws1.Range("A4").Value = "Y"
$ws1.Range("B4").Value = "Medical, Health & Beauty"
$ws1.Range("C4").Value = "medical-health-beauty"

# --- Sheet1: add a second report group underneath ------------------------
$ws1.Range("A7").Value = "CheckTitleOfHomePageLink2"

$ws1.Range("A8").Value = "Runmode"
$ws1.Range("B8").Value = "Link"
$ws1.Range("C8").Value = "Title"

$ws1.Range("A9").Value = "Y"
$ws1.Range("B9").Value = "Machinery & Industry"
$ws1.Range("A10").Value = "Y"
$ws1.Range("B10").Value = "Automotive"
$ws1.Range("C9").Value = "machinery_industrial_plant_hardware_parts_tools"
$ws1.Range("C10").Value = "automotive_vehicles_transportation_parts"

# --- Add Sheet2 with a CheckStaticLinks report group ---------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

$ws2.Range("A1").Value = "CheckStaticLinks"
$ws2.Range("B3").Value = "About Us"
$ws2.Range("A2").Value = "RunMode"
$ws2.Range("B2").Value = "Expected Data"
$ws2.Range("B4").Value = "Help"
$ws2.Range("B5").Value = "News"
$ws2.Range("A3").Value = "Y"
$ws2.Range("A4").Value = "Y"
$ws2.Range("A5").Value = "Y"

# --- Selections / active sheet -------------------------------------------
$ws1.Range("A11").Select()
$ws2.Select()
